$wb = $excel.ActiveWorkbook

# Sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1464.6296
$ws.Range("I129").Value = 695
$ws.Range("J129").Value = 1494.2307
$ws.Range("K129").Value = 2085
$ws.Range("L129").Value = 4482.6921
$ws.Range("M129").Value = 2915
$ws.Range("N129").Value = -14482.6921

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 13554.863
$ws.Range("I137").Value = 960.06665
$ws.Range("J137").Value = 40543.715
$ws.Range("K137").Value = 2880.19995
$ws.Range("L137").Value = 121631.145
$ws.Range("M137").Value = -330.1999500000002
$ws.Range("N137").Value = -126731.145

# Sheet ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1500
$ws.Range("I141").Value = 1312.5
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 3937.5
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 1242.5
$ws.Range("N141").Value = -16360

# Sheet ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 362.3846
$ws.Range("I5").Value = 354.625
$ws.Range("J5").Value = 374.8
$ws.Range("K5").Value = 354.625
$ws.Range("L5").Value = 374.8
$ws.Range("M5").Value = -242.625
$ws.Range("N5").Value = -598.8

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14778.282
$ws.Range("I32").Value = 13995.578
$ws.Range("K32").Value = 13995.578
$ws.Range("M32").Value = -13708.578

# Sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2320.9167
$ws.Range("I61").Value = 2204.4348
$ws.Range("K61").Value = 2204.4348
$ws.Range("M61").Value = -1992.4348

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4704.5
$ws.Range("I74").Value = 1200.5454
$ws.Range("J74").Value = 17552.334
$ws.Range("K74").Value = 1200.5454
$ws.Range("L74").Value = 17552.334
$ws.Range("M74").Value = -326.5454
$ws.Range("N74").Value = -19300.334

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4704.5
$ws.Range("I77").Value = 1200.5454
$ws.Range("J77").Value = 17552.334
$ws.Range("K77").Value = 6002.727
$ws.Range("L77").Value = 87761.67
$ws.Range("M77").Value = -1634.727
$ws.Range("N77").Value = -96497.67

# Sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2320.9167
$ws.Range("I136").Value = 2204.4348
$ws.Range("K136").Value = 6613.3044
$ws.Range("M136").Value = -4063.3044

# Sheet BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 362.3846
$ws.Range("I4").Value = 354.625
$ws.Range("J4").Value = 374.8
$ws.Range("K4").Value = 354.625
$ws.Range("L4").Value = 374.8
$ws.Range("M4").Value = -239.625
$ws.Range("N4").Value = -604.8

# Sheet BSM row 62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 32565
$ws.Range("J62").Value = 32565
$ws.Range("L62").Value = 32565
$ws.Range("N62").Value = -33937

# Sheet BSM row 65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 32565
$ws.Range("J65").Value = 32565
$ws.Range("L65").Value = 97695
$ws.Range("N65").Value = -104559

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 866.76666
$ws.Range("I94").Value = 640.9583
$ws.Range("J94").Value = 1770
$ws.Range("K94").Value = 640.9583
$ws.Range("L94").Value = 1770
$ws.Range("M94").Value = -189.9583
$ws.Range("N94").Value = -2672

# Sheet CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7750000
$ws.Range("I6").Value = 7428571.5
$ws.Range("K6").Value = 7428571.5
$ws.Range("M6").Value = -7428458.5

# Sheet CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.2
$ws.Range("I7").Value = 50.666668
$ws.Range("J7").Value = 533.3333
$ws.Range("K7").Value = 50.666668
$ws.Range("L7").Value = 533.3333
$ws.Range("M7").Value = 62.333332
$ws.Range("N7").Value = -759.3333

# Sheet CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 16168.857
$ws.Range("I16").Value = 20536
$ws.Range("J16").Value = 5251
$ws.Range("K16").Value = 20536
$ws.Range("L16").Value = 5251
$ws.Range("M16").Value = -20249
$ws.Range("N16").Value = -5825

# Sheet CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 22666.666
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 30000
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = -7572
$ws.Range("N41").Value = -30856

# Sheet CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9243
$ws.Range("J50").Value = 9243
$ws.Range("L50").Value = 9243
$ws.Range("N50").Value = -10493

# Sheet CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9251.6
$ws.Range("J51").Value = 9251.6
$ws.Range("L51").Value = 9251.6
$ws.Range("N51").Value = -10723.6

# Sheet CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 16009.75
$ws.Range("J59").Value = 16009.75
$ws.Range("L59").Value = 16009.75
$ws.Range("N59").Value = -18299.75

# Sheet CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 7637.125
$ws.Range("I60").Value = 6666.3335
$ws.Range("J60").Value = 8219.6
$ws.Range("K60").Value = 6666.3335
$ws.Range("L60").Value = 8219.6
$ws.Range("M60").Value = -6155.3335
$ws.Range("N60").Value = -9241.6

# Sheet CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9251.6
$ws.Range("J61").Value = 9251.6
$ws.Range("L61").Value = 9251.6
$ws.Range("N61").Value = -9947.6

# Sheet CRP row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Sheet CRP row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Sheet CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17242.428
$ws.Range("J68").Value = 17242.428
$ws.Range("L68").Value = 17242.428
$ws.Range("N68").Value = -18740.428

# Sheet CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 17242.428
$ws.Range("J71").Value = 17242.428
$ws.Range("L71").Value = 51727.284
$ws.Range("N71").Value = -59215.284

# Sheet CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 13796.857
$ws.Range("J74").Value = 13796.857
$ws.Range("L74").Value = 13796.857
$ws.Range("N74").Value = -15544.857

# Sheet CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 13796.857
$ws.Range("J77").Value = 13796.857
$ws.Range("L77").Value = 41390.571
$ws.Range("N77").Value = -50126.571

# Sheet CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 16168.857
$ws.Range("I113").Value = 20536
$ws.Range("J113").Value = 5251
$ws.Range("K113").Value = 20536
$ws.Range("L113").Value = 5251
$ws.Range("M113").Value = -18366
$ws.Range("N113").Value = -9591

# Sheet CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10622
$ws.Range("M69").ClearContents()

# Sheet CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35112
$ws.Range("M72").ClearContents()

# Sheet GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13374.2
$ws.Range("I70").Value = 21245.273
$ws.Range("J70").Value = 3754
$ws.Range("K70").Value = 21245.273
$ws.Range("L70").Value = 3754
$ws.Range("M70").Value = -20975.273
$ws.Range("N70").Value = -4294

# Sheet GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13374.2
$ws.Range("I73").Value = 21245.273
$ws.Range("J73").Value = 3754
$ws.Range("K73").Value = 21245.273
$ws.Range("L73").Value = 3754
$ws.Range("M73").Value = -20309.273
$ws.Range("N73").Value = -5626
